# Updates betting-odds values on Sheet1 to match the 2025-06-04 FlashScore refresh.
# Only numeric odds cells change; row/column layout and headers are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("K2").Value = 8
$ws.Range("N2").Value = 2.2
$ws.Range("O2").Value = 1.65

# Row 5
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 4.2

# Row 7
$ws.Range("G7").Value = 1.47
$ws.Range("H7").Value = 4.3
$ws.Range("I7").Value = 5.7
$ws.Range("N7").Value = 1.57
$ws.Range("O7").Value = 2.12
$ws.Range("R7").Value = 1.57
$ws.Range("S7").Value = 2.1
$ws.Range("T7").Value = 9.5
$ws.Range("U7").Value = 8.75
$ws.Range("V7").Value = 8.25
$ws.Range("W7").Value = 11.25
$ws.Range("X7").Value = 10.75
$ws.Range("Y7").Value = 19
$ws.Range("Z7").Value = 16.5
$ws.Range("AA7").Value = 9
$ws.Range("AB7").Value = 14.5
$ws.Range("AC7").Value = 50
$ws.Range("AD7").Value = 300
$ws.Range("AE7").Value = 21
$ws.Range("AG7").Value = 18
$ws.Range("AH7").Value = 110
$ws.Range("AI7").Value = 50
$ws.Range("AJ7").Value = 45

# Row 8
$ws.Range("G8").Value = 2.35
$ws.Range("H8").Value = 2.92
$ws.Range("I8").Value = 3.1
$ws.Range("L8").Value = 1.35
$ws.Range("M8").Value = 2.7
$ws.Range("N8").Value = 2.02
$ws.Range("O8").Value = 1.62
$ws.Range("P8").Value = 1.42
$ws.Range("Q8").Value = 2.47
$ws.Range("T8").Value = 7.4
$ws.Range("U8").Value = 11.5
$ws.Range("V8").Value = 9
$ws.Range("W8").Value = 25
$ws.Range("X8").Value = 20
$ws.Range("Y8").Value = 30
$ws.Range("Z8").Value = 7.9
$ws.Range("AA8").Value = 5.7
$ws.Range("AB8").Value = 13.5
$ws.Range("AC8").Value = 65
$ws.Range("AE8").Value = 8.5
$ws.Range("AF8").Value = 16
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 45
$ws.Range("AI8").Value = 29
$ws.Range("AJ8").Value = 37

# Row 10
$ws.Range("H10").Value = 3.7
$ws.Range("I10").Value = 2.65
$ws.Range("N10").Value = 1.57
$ws.Range("O10").Value = 2.12
$ws.Range("R10").Value = 1.52
$ws.Range("S10").Value = 2.22
$ws.Range("T10").Value = 10.5
$ws.Range("U10").Value = 13
$ws.Range("X10").Value = 16.5
$ws.Range("Y10").Value = 22
$ws.Range("Z10").Value = 15
$ws.Range("AA10").Value = 7.5
$ws.Range("AB10").Value = 12.5
$ws.Range("AC10").Value = 45
$ws.Range("AD10").Value = 250
$ws.Range("AE10").Value = 11.75
$ws.Range("AF10").Value = 15.5
$ws.Range("AG10").Value = 10
$ws.Range("AI10").Value = 19.5
$ws.Range("AJ10").Value = 24

# Row 12
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 3.75
$ws.Range("M12").Value = 2.35
$ws.Range("N12").Value = 2.35
$ws.Range("O12").Value = 1.47
$ws.Range("P12").Value = 1.5
$ws.Range("Q12").Value = 2.25
$ws.Range("R12").Value = 2.07
$ws.Range("S12").Value = 1.6
$ws.Range("T12").Value = 5.5
$ws.Range("U12").Value = 8.25
$ws.Range("V12").Value = 9.25
$ws.Range("W12").Value = 17.5
$ws.Range("X12").Value = 20
$ws.Range("Y12").Value = 40
$ws.Range("Z12").Value = 6.7
$ws.Range("AA12").Value = 6.2
$ws.Range("AB12").Value = 19
$ws.Range("AE12").Value = 8.25
$ws.Range("AF12").Value = 18.5
$ws.Range("AG12").Value = 13.5
$ws.Range("AH12").Value = 60
$ws.Range("AI12").Value = 45
$ws.Range("AJ12").Value = 60
